$wb = $excel.ActiveWorkbook
$ref = $wb.Worksheets.Item("2021-Q4")

# --- Step 1: insert the new "2022-Q1" sheet before "总计" ---
$total = $wb.Worksheets.Item("总计")
$ws = $wb.Worksheets.Add($total)
$ws.Name = "2022-Q1"

# Headers (row 1)
$ws.Cells.Item(1, 2).Value = "基金代码"
$ws.Cells.Item(1, 3).Value = "基金名称"
$ws.Cells.Item(1, 4).Value = "基金规模"
$ws.Cells.Item(1, 5).Value = "股票总仓位"
$ws.Cells.Item(1, 6).Value = "仓位占比"
$ws.Cells.Item(1, 7).Value = "持有市值(亿元)"
$ws.Cells.Item(1, 8).Value = "仓位排名"

# Copy header style (bold/border/center) from the reference sheet
$ref.Range("B1:H1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)

# Row 2: fund 001027
$rng = $ws.Range("B2:G2")
$rng.NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "001027"
$ws.Cells.Item(2, 3).Value = "前海开源中证大农业指数增强"
$ws.Cells.Item(2, 4).Value = "6.63"
$ws.Cells.Item(2, 5).Value = "92.76"
$ws.Cells.Item(2, 6).Value = "6.48"
$ws.Cells.Item(2, 7).Value = "0.4296"
$ws.Cells.Item(2, 8).Value = 2
$rng.ClearFormats()

# Row 3: fund 005347
$rng = $ws.Range("B3:G3")
$rng.NumberFormat = "@"
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = "005347"
$ws.Cells.Item(3, 3).Value = "诺德量化优选6个月持有期混合"
$ws.Cells.Item(3, 4).Value = "2.60"
$ws.Cells.Item(3, 5).Value = "93.66"
$ws.Cells.Item(3, 6).Value = "2.95"
$ws.Cells.Item(3, 7).Value = "0.0767"
$ws.Cells.Item(3, 8).Value = 7
$rng.ClearFormats()

# Row 4: fund 006440
$rng = $ws.Range("B4:G4")
$rng.NumberFormat = "@"
$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = "006440"
$ws.Cells.Item(4, 3).Value = "中信建投中证500指数增强A"
$ws.Cells.Item(4, 4).Value = "5.78"
$ws.Cells.Item(4, 5).Value = "94.71"
$ws.Cells.Item(4, 6).Value = "0.99"
$ws.Cells.Item(4, 7).Value = "0.0572"
$ws.Cells.Item(4, 8).Value = 7
$rng.ClearFormats()

# Row 5: fund 006267
$rng = $ws.Range("B5:G5")
$rng.NumberFormat = "@"
$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = "006267"
$ws.Cells.Item(5, 3).Value = "诺德量化核心灵活配置混合A"
$ws.Cells.Item(5, 4).Value = "1.84"
$ws.Cells.Item(5, 5).Value = "93.91"
$ws.Cells.Item(5, 6).Value = "2.98"
$ws.Cells.Item(5, 7).Value = "0.0548"
$ws.Cells.Item(5, 8).Value = 7
$rng.ClearFormats()

# Row 6: fund 006441
$rng = $ws.Range("B6:G6")
$rng.NumberFormat = "@"
$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = "006441"
$ws.Cells.Item(6, 3).Value = "中信建投中证500指数增强C"
$ws.Cells.Item(6, 4).Value = "3.11"
$ws.Cells.Item(6, 5).Value = "94.71"
$ws.Cells.Item(6, 6).Value = "0.99"
$ws.Cells.Item(6, 7).Value = "0.0308"
$ws.Cells.Item(6, 8).Value = 7
$rng.ClearFormats()

# Row 7: fund 006268
$rng = $ws.Range("B7:G7")
$rng.NumberFormat = "@"
$ws.Cells.Item(7, 1).Value = 5
$ws.Cells.Item(7, 2).Value = "006268"
$ws.Cells.Item(7, 3).Value = "诺德量化核心灵活配置混合C"
$ws.Cells.Item(7, 4).Value = "0.50"
$ws.Cells.Item(7, 5).Value = "93.91"
$ws.Cells.Item(7, 6).Value = "2.98"
$ws.Cells.Item(7, 7).Value = "0.0149"
$ws.Cells.Item(7, 8).Value = 7
$rng.ClearFormats()

# Row 8: fund 004192
$rng = $ws.Range("B8:G8")
$rng.NumberFormat = "@"
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = "004192"
$ws.Cells.Item(8, 3).Value = "招商中证500指数增强A"
$ws.Cells.Item(8, 4).Value = "0.96"
$ws.Cells.Item(8, 5).Value = "94.32"
$ws.Cells.Item(8, 6).Value = "1.06"
$ws.Cells.Item(8, 7).Value = "0.0102"
$ws.Cells.Item(8, 8).Value = 6
$rng.ClearFormats()

# Row 9: fund 004193
$rng = $ws.Range("B9:G9")
$rng.NumberFormat = "@"
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "004193"
$ws.Cells.Item(9, 3).Value = "招商中证500指数增强C"
$ws.Cells.Item(9, 4).Value = "0.42"
$ws.Cells.Item(9, 5).Value = "94.32"
$ws.Cells.Item(9, 6).Value = "1.06"
$ws.Cells.Item(9, 7).Value = "0.0045"
$ws.Cells.Item(9, 8).Value = 6
$rng.ClearFormats()

# Row 10: fund 002453
$rng = $ws.Range("B10:G10")
$rng.NumberFormat = "@"
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = "002453"
$ws.Cells.Item(10, 3).Value = "九泰久稳灵活配置混合A"
$ws.Cells.Item(10, 4).Value = "0.09"
$ws.Cells.Item(10, 5).Value = "94.85"
$ws.Cells.Item(10, 6).Value = "2.25"
$ws.Cells.Item(10, 7).Value = "0.0020"
$ws.Cells.Item(10, 8).Value = 6
$rng.ClearFormats()

# Row 11: fund 002454
$rng = $ws.Range("B11:G11")
$rng.NumberFormat = "@"
$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = "002454"
$ws.Cells.Item(11, 3).Value = "九泰久稳灵活配置混合C"
$ws.Cells.Item(11, 4).Value = "0.04"
$ws.Cells.Item(11, 5).Value = "94.85"
$ws.Cells.Item(11, 6).Value = "2.25"
$ws.Cells.Item(11, 7).Value = "0.0009"
$ws.Cells.Item(11, 8).Value = 6
$rng.ClearFormats()

# Apply column-A style (bold/border/center) to match the other sheets
$ref.Range("A2").Copy()
$ws.Range("A2:A11").PasteSpecial(-4122)

# --- Step 2: insert a new top row into "总计" for the 2022-Q1 summary ---
$tot = $wb.Worksheets.Item("总计")
$tot.Rows(2).Insert()
$clr = $tot.Range("B2:D2")
$clr.ClearFormats()
$tot.Cells.Item(2, 1).Value = 0
$tot.Cells.Item(2, 2).Value = "2022-Q1"
$tot.Cells.Item(2, 3).Value = 10
$tot.Cells.Item(2, 4).Value = 0.68
$tot.Range("A3").Copy()
$tot.Range("A2").PasteSpecial(-4122)

Write-Host "Edit complete"